$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157, pushing existing rows 157:178 down to 158:179.
$ws.Rows("157:157").Insert()

# The new row 157 is a new weekly record; copy the non-changing columns from the
# row that is now directly below it (old row 157, now at 158) and set the new date.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(157, $c).Value2 = $ws.Cells.Item(158, $c).Value2
}

$ws.Cells.Item(157, 4).Value2 = 44474
$ws.Cells.Item(157, 10).Value2 = 500
$ws.Cells.Item(157, 11).Value2 = 1800
$ws.Cells.Item(157, 12).Value2 = 1800
$ws.Cells.Item(157, 13).Value2 = 1800
$ws.Cells.Item(157, 16).Value2 = 450
